# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new F value, applied identically to both sheets.
$updates = @{
    2  = 196
    3  = 3309
    4  = 242
    6  = 208
    7  = 1721
    8  = 1652
    9  = 472
    11 = 234
    17 = 241
    21 = 20
    22 = 61
    23 = 19
    25 = 392
    26 = 248
    27 = 111
    29 = 14
    31 = 372
    32 = 2256
    35 = 477
    36 = 550
    39 = 232
    42 = 536
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
